$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting all existing columns right
$ws.Columns.Item(1).Insert()

# Set header text for the newly inserted column
$ws.Range("A1").Value = "Preparer"

# Move selection to A2, mirroring the state after typing the header and pressing Enter
$ws.Range("A2").Select()
